$d = $word.ActiveDocument

# 1. Merge "What I plan to do before the next progress " + "report:" (which
#    were split across two runs with a _GoBack bookmark sandwiched between
#    them) into a single plain run of text.
$d.Content.Find.Execute("What I plan to do before the next progress report:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "What I plan to do before the next progress report:", 2) | Out-Null

# 2. Fill in the Team Progress score of "2": the run of 42 underscores
#    following "Team Progress (0 -- 3):  " is split into "__" + an
#    underlined "2" + 39 remaining underscores. The new "2" run is wrapped
#    in the (relocated) _GoBack bookmark.
$rng = $d.Content
$rng.Find.Execute("Team Progress (0 -- 3):  ") | Out-Null
$usStart = $rng.End

$two = $d.Range($usStart + 2, $usStart + 3)
$two.Text = "2"
$two.Font.Underline = 1
$d.Bookmarks.Add("_GoBack", $two) | Out-Null
